$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Reword the existing row 9 story text and switch its story-text cell
# from the Calibri wrap style to the Arial wrap style used by the rest of
# the table (copy format from B2, which already carries that style). ---
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "As a student looking for an apartment, I want to sort the apartment offerings by price so that I can afford easily."

# --- Fill in the new user-story rows (10-27), copying the existing row's
# formatting (A/B/C/D/F from row 2, E - status - from row 6) so the cell
# styles line up with the rest of the backlog table. ---
# Row 10
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A10:F10").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "As a landlord, I want to advertise one of my apartments on the apartment section so that students can find and rent my apartment."
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "M"
$ws.Range("E10").Value = "W"
$ws.Range("F10").Value = 3

# Row 11
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A11:F11").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "As a student looking to sell my textbook, I want to list last year's textbooks for sale so that I can get some money for them."
$ws.Range("C11").Value = 4
$ws.Range("D11").Value = "M"
$ws.Range("E11").Value = "W"
$ws.Range("F11").Value = 4

# Row 12
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A12:F12").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "As a student looking to sell my textbook, I want a list of how much users paid for this textbook last semester so I can make money."
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = "M"
$ws.Range("E12").Value = "W"
$ws.Range("F12").Value = 5

# Row 13
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A13:F13").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "As a student selling items, I want to edit the tags on my listing so that they can get better exposure."
$ws.Range("C13").Value = 6
$ws.Range("D13").Value = "M"
$ws.Range("E13").Value = "W"
$ws.Range("F13").Value = 6

# Row 14
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A14:F14").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "As a student looking for furniture, I want to use CSwap and navigate to the furniture section so that I can find furniture for my apartment."
$ws.Range("C14").Value = 7
$ws.Range("D14").Value = "M"
$ws.Range("E14").Value = "W"
$ws.Range("F14").Value = 7

# Row 15
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A15:F15").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "As a student looking for furniture, I want to sort the furniture section so that I can find only what I need."
$ws.Range("C15").Value = 8
$ws.Range("D15").Value = "M"
$ws.Range("E15").Value = "W"
$ws.Range("F15").Value = 8

# Row 16
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A16:F16").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "As a seller, I want to list my calculator for sale under the electronics section so that other students can buy it."
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = "M"
$ws.Range("E16").Value = "W"
$ws.Range("F16").Value = 9

# Row 17
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A17:F17").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "As a user, I want a profile so that I can see what the person I am selling to looks like."
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = "M"
$ws.Range("E17").Value = "W"
$ws.Range("F17").Value = 10

# Row 18
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A18:F18").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "As a seller, I want to add new photos to my listing so that my listing is more attractive."
$ws.Range("C18").Value = 11
$ws.Range("D18").Value = "M"
$ws.Range("E18").Value = "W"
$ws.Range("F18").Value = 11

# Row 19
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A19:F19").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "As a user, I want to delete my account so that I am no longer active on the app."
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = "M"
$ws.Range("E19").Value = "W"
$ws.Range("F19").Value = 12

# Row 20
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A20:F20").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "As a seller, I want to change the description of my listing so that It’s available for both pick up and for delivery."
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = "M"
$ws.Range("E20").Value = "W"
$ws.Range("F20").Value = 13

# Row 21
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A21:F21").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E21").PasteSpecial(-4122) | Out-Null
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "As a landlord, I want to be able to list my apartment under multiple colleges so that my apartment will be seen by people in the area."
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = "M"
$ws.Range("E21").Value = "W"
$ws.Range("F21").Value = 14

# Row 22
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A22:F22").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "As a seller, I want to be able to see all the items I am currently selling so that I know from organizational purposes."
$ws.Range("C22").Value = 15
$ws.Range("D22").Value = "M"
$ws.Range("E22").Value = "W"
$ws.Range("F22").Value = 15
$ws.Range("B22").VerticalAlignment = -4160

# Row 23
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A23:F23").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "As a seller, I want to see how long my listing has been on the website so that I can determine if I need to make a change to the listing."
$ws.Range("C23").Value = 16
$ws.Range("D23").Value = "M"
$ws.Range("E23").Value = "W"
$ws.Range("F23").Value = 16
$ws.Range("B23").VerticalAlignment = -4160

# Row 24
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A24:F24").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E24").PasteSpecial(-4122) | Out-Null
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "As a user, I want to be able to give ratings to a seller so that I can help provide useful feedback on the buying process."
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = "M"
$ws.Range("E24").Value = "W"
$ws.Range("F24").Value = 17
$ws.Range("B24").VerticalAlignment = -4160

# Row 25
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A25:F25").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "As a user, I want to be able to see my recently viewed listings so that I can go back to view listings I previously visited."
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = "M"
$ws.Range("E25").Value = "W"
$ws.Range("F25").Value = 18
$ws.Range("B25").VerticalAlignment = -4108

# Row 26
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A26:F26").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "As a student, I want to be able to have a book condition tag so that I can sell my high-quality books for more money."
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = "M"
$ws.Range("E26").Value = "W"
$ws.Range("F26").Value = 19
$ws.Range("B26").VerticalAlignment = -4160

# Row 27
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A27:F27").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "As a student user, I want to filter out listings above a price range so I don't only see listings I can't afford."
$ws.Range("C27").Value = 20
$ws.Range("D27").Value = "M"
$ws.Range("E27").Value = "W"
$ws.Range("F27").Value = 20

# --- Append the blank template rows below the table. ---
# Row 28 keeps the Arial (row-2-style) story-text formatting.
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A28:F28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Clear() | Out-Null

# Rows 29-42 revert to the original Calibri-wrap blank-row style
# (copied from the sheet's pre-existing blank row 10/B10 pattern).
# Column E is intentionally left untouched (blank / unstyled) on
# every one of these rows, matching the source workbook.
# Row 29
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C29:D29").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null

# Row 30
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C30:D30").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null

# Row 31
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C31:D31").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null

# Row 32
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A32").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B32").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C32:D32").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F32").PasteSpecial(-4122) | Out-Null

# Row 33
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A33").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C33:D33").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null

# Row 34
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A34").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C34:D34").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F34").PasteSpecial(-4122) | Out-Null

# Row 35
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A35").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B35").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C35:D35").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F35").PasteSpecial(-4122) | Out-Null

# Row 36
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A36").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B36").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C36:D36").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F36").PasteSpecial(-4122) | Out-Null

# Row 37
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C37:D37").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F37").PasteSpecial(-4122) | Out-Null

# Row 38
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B38").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C38:D38").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F38").PasteSpecial(-4122) | Out-Null

# Row 39
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A39").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B39").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C39:D39").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F39").PasteSpecial(-4122) | Out-Null

# Row 40
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A40").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B40").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C40:D40").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F40").PasteSpecial(-4122) | Out-Null

# Row 41
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A41").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B41").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C41:D41").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F41").PasteSpecial(-4122) | Out-Null

# Row 42
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A42").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B42").PasteSpecial(-4122) | Out-Null
$ws.Range("C2:D2").Copy() | Out-Null
$ws.Range("C42:D42").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F42").PasteSpecial(-4122) | Out-Null

Write-Host "Done applying Product Backlog edits"
